$wb = $excel.ActiveWorkbook

# Remove the empty "Sheet" worksheet so only "Log" remains
$excel.DisplayAlerts = $false
$wb.Worksheets.Item("Sheet").Delete()

$ws = $wb.Worksheets.Item("Log")

# Style the header row (A1:B1): bold font, thin box border, centered + top aligned.
# Apply per-cell (rather than on the A1:B1 range in one shot) so both cells
# converge on the exact same cell style.
foreach ($addr in @("A1", "B1")) {
    $cell = $ws.Range($addr)
    $cell.Font.Bold = $true
    $cell.Borders.LineStyle = 1   # xlContinuous
    $cell.Borders.Weight = 2      # xlThin
    $cell.HorizontalAlignment = -4108  # xlCenter
    $cell.VerticalAlignment = -4160    # xlTop
}

# Append the new log entry row
$ws.Range("A2").Value = "2025-05-28 11:42:34"
$ws.Range("B2").Value = ""
